$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / URL / label cells (safe to assign directly as strings) ---
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E9").Value = "8OneONEBestin24h"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

# --- Numeric-looking price cells: force text storage (preserve exact
#     decimal formatting such as trailing/leading zeros) by using the
#     leading-apostrophe text-entry convention, then reset the style back
#     to Normal so no stray number-format style is left on the cell. ---
$ws.Range("D2").Value = "'244.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("D4").Value = "'5.388"
$ws.Range("D4").Style = "Normal"
$ws.Range("D6").Value = "'3.393"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Value = "'0.9276"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.01126"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1433"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07503"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03488"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03049"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09431"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'4.016"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001589"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04806"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.005611"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.004170"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.0009917"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'3.666"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'6.427"
$ws.Range("D22").Style = "Normal"
$ws.Range("D26").Value = "'0.00007000"
$ws.Range("D26").Style = "Normal"
$ws.Range("D41").Value = "'0.006380"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Value = "'0.002900"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.005900"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005245"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'1.0000"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.002330"
$ws.Range("D48").Style = "Normal"
